$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values (Strike# renamed to K). Update rows 2-9 per regen.
$values = @{
    2 = 2
    3 = 6
    4 = 8
    5 = 4
    6 = 9
    7 = 5
    8 = 2
    9 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
